# BOT; UPDATE DATA
# Adds the 2020-05-09 (Excel serial 43960) daily row to the "all", "kobe"
# and "other" sheets, ahead of the trailing footnote row on each sheet,
# and refreshes the "kobe" sheet's previous-day (43959) cumulative totals
# for columns D/E.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": insert new row 32 (date 43960) before the footnote row,
# which shifts the footnote down to row 33.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()
$wsAll.Rows.Item(32).Insert()

$wsAll.Range("A32").Value = 43960
$wsAll.Range("B32").Value = 276
$wsAll.Range("C32").Value = 271
$wsAll.Range("D32").Value = 80
$wsAll.Range("E32").Value = 70
$wsAll.Range("F32").Value = 10
$wsAll.Range("G32").Value = 8
$wsAll.Range("H32").Value = 183

$excel.ActiveWindow.FreezePanes = $false
$wsAll.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsAll.Range("A32").Select()

# ---------------------------------------------------------------------
# Sheet "kobe": update the 43959 row's D/E cumulative totals, then
# insert new row 87 (date 43960) before the footnote row, which shifts
# the footnote down to row 88.
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate()

$wsKobe.Range("D86").Value = 3
$wsKobe.Range("E86").Value = 276

$wsKobe.Rows.Item(87).Insert()

$wsKobe.Range("A87").Value = 43960
$wsKobe.Range("B87").Value = 0
$wsKobe.Range("C87").Value = 2480
$wsKobe.Range("D87").Value = 0
$wsKobe.Range("E87").Value = 276
$wsKobe.Range("F87").Value = 75
$wsKobe.Range("G87").Value = 66
$wsKobe.Range("H87").Value = 9
$wsKobe.Range("I87").Value = 8
$wsKobe.Range("J87").Value = 174

$wsKobe.Range("A87").Select()

# ---------------------------------------------------------------------
# Sheet "other": insert new row 62 (date 43960) before the footnote
# row, which shifts the footnote down to row 63.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()
$wsOther.Rows.Item(62).Insert()

$wsOther.Range("A62").Value = 43960
$wsOther.Range("B62").Value = 0
$wsOther.Range("C62").Value = 14
$wsOther.Range("D62").Value = 5
$wsOther.Range("E62").Value = 4
$wsOther.Range("F62").Value = 1
$wsOther.Range("G62").Value = 0
$wsOther.Range("H62").Value = 9

$wsOther.Range("A62").Select()

$wsAll.Activate()

Write-Host "Applied daily update for 2020-05-09 (serial 43960) to all/kobe/other sheets."
